$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = $true
$ws.Range("A9").Value = $false

$ws.Range("A10").Select() | Out-Null
